$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A formatting (border/bold/center style already used by A2:A3) down to A4:A5
$ws.Range("A3").Copy($ws.Range("A4:A5"))

# Row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 35

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 27

# Row 4 (new)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 23

# Row 5 (new)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 14
